# Generate Report for Handback
# The handback transform failed for a76a6d29-a487-4709-b493-3cd9d17f4647.md
# because the handback file name didn't match the expected handoff file name.
# Update the status to "Handback transform failed" and record the per-language
# error detail, widening the Error Detail column so the message is readable.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"
$zhError = "Handback file name: e02sqrhq.3up is different with handoff file name: a76a6d29-a487-4709-b493-3cd9d17f4647.dc3ef1ec36b78102e6e28f6f2c6b1acd59cf3f44.zh-cn."
$deError  = "Handback file name: e02sqrhq.3up is different with handoff file name: a76a6d29-a487-4709-b493-3cd9d17f4647.dc3ef1ec36b78102e6e28f6f2c6b1acd59cf3f44.de-de."

# --- Overview sheet: reflect the new status for the zh-cn / de-de columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E7").Value = $newStatus
$wsOverview.Range("F7").Value = $newStatus

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C7").Value = $newStatus
$wsZh.Range("P7").Value = $zhError
$wsZh.Range("P1").ColumnWidth = 39.17

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C7").Value = $newStatus
$wsDe.Range("P7").Value = $deError
$wsDe.Range("P1").ColumnWidth = 39.17
